# Auto-generated edit script applying numeric updates to Aegis_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 694.6818
$ws.Range("I28").Value = 400.92856
$ws.Range("J28").Value = 1208.75
$ws.Range("K28").Value = 400.92856
$ws.Range("L28").Value = 1208.75
$ws.Range("M28").Value = 84.07144
$ws.Range("N28").Value = -2178.75

$ws.Range("H40").Value = 1958
$ws.Range("I40").Value = 1344.5454
$ws.Range("J40").Value = 2477.077
$ws.Range("K40").Value = 1344.5454
$ws.Range("L40").Value = 2477.077
$ws.Range("M40").Value = -1169.5454
$ws.Range("N40").Value = -2827.077

$ws.Range("H76").Value = 4683.5557
$ws.Range("I76").Value = 4644
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4644
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4329
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 4683.5557
$ws.Range("I79").Value = 4644
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4644
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3552
$ws.Range("N79").Value = -7184

$ws.Range("H103").Value = 2697.5557
$ws.Range("J103").Value = 1932.5714
$ws.Range("L103").Value = 5797.7142
$ws.Range("N103").Value = -6969.7142

$ws.Range("H132").Value = 5213630.5
$ws.Range("I132").Value = 5687072
$ws.Range("J132").Value = 5777.75
$ws.Range("K132").Value = 17061216
$ws.Range("L132").Value = 17333.25
$ws.Range("M132").Value = -17058686
$ws.Range("N132").Value = -22393.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28812.305
$ws.Range("I2").Value = 1161.7037
$ws.Range("J2").Value = 111764.11
$ws.Range("K2").Value = 1161.7037
$ws.Range("L2").Value = 111764.11
$ws.Range("M2").Value = -1048.7037
$ws.Range("N2").Value = -111990.11

$ws.Range("H32").Value = 24313.53
$ws.Range("I32").Value = 5787.831
$ws.Range("K32").Value = 5787.831
$ws.Range("M32").Value = -5500.831

$ws.Range("H63").Value = 2505
$ws.Range("I63").Value = 2219.875
$ws.Range("J63").Value = 2961.2
$ws.Range("K63").Value = 2219.875
$ws.Range("L63").Value = 2961.2
$ws.Range("M63").Value = -1533.875
$ws.Range("N63").Value = -4333.2

$ws.Range("H66").Value = 2505
$ws.Range("I66").Value = 2219.875
$ws.Range("J66").Value = 2961.2
$ws.Range("K66").Value = 11099.375
$ws.Range("L66").Value = 14806
$ws.Range("M66").Value = -7667.375
$ws.Range("N66").Value = -21670

$ws.Range("H97").Value = 26792.95
$ws.Range("I97").Value = 37764.184
$ws.Range("J97").Value = 2107.6667
$ws.Range("K97").Value = 37764.184
$ws.Range("L97").Value = 2107.6667
$ws.Range("M97").Value = -37268.184
$ws.Range("N97").Value = -3099.6667

$ws.Range("H116").Value = 28812.305
$ws.Range("I116").Value = 1161.7037
$ws.Range("J116").Value = 111764.11
$ws.Range("K116").Value = 1161.7037
$ws.Range("L116").Value = 111764.11
$ws.Range("M116").Value = 1132.2963
$ws.Range("N116").Value = -116352.11

$ws.Range("H132").Value = 10981.917
$ws.Range("I132").Value = 13202.195
$ws.Range("J132").Value = 3686.7144
$ws.Range("K132").Value = 39606.585
$ws.Range("L132").Value = 11060.1432
$ws.Range("M132").Value = -37076.585
$ws.Range("N132").Value = -16120.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28812.305
$ws.Range("I3").Value = 1161.7037
$ws.Range("J3").Value = 111764.11
$ws.Range("K3").Value = 1161.7037
$ws.Range("L3").Value = 111764.11
$ws.Range("M3").Value = -1047.7037
$ws.Range("N3").Value = -111992.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1599.1613
$ws.Range("I58").Value = 1348.4615
$ws.Range("K58").Value = 1348.4615
$ws.Range("M58").Value = -1145.4615

$ws.Range("H99").Value = 8689.3125
$ws.Range("I99").Value = 2120
$ws.Range("K99").Value = 2120
$ws.Range("M99").Value = -622

$ws.Range("H126").Value = 8689.3125
$ws.Range("I126").Value = 2120
$ws.Range("K126").Value = 6360
$ws.Range("M126").Value = -3890

$ws.Range("H136").Value = 1599.1613
$ws.Range("I136").Value = 1348.4615
$ws.Range("K136").Value = 4045.3845
$ws.Range("M136").Value = -1495.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1839.1852
$ws.Range("I132").Value = 810.4
$ws.Range("J132").Value = 2444.353
$ws.Range("K132").Value = 7293.599999999999
$ws.Range("L132").Value = 21999.177
$ws.Range("M132").Value = -4763.599999999999
$ws.Range("N132").Value = -27059.177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1522.7354
$ws.Range("I122").Value = 1271.5
$ws.Range("K122").Value = 3814.5
$ws.Range("M122").Value = -1364.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 59953.668
$ws.Range("J33").Value = 85021
$ws.Range("L33").Value = 85021
$ws.Range("N33").Value = -85521

$ws.Range("H36").Value = 59953.668
$ws.Range("J36").Value = 85021
$ws.Range("L36").Value = 85021
$ws.Range("N36").Value = -85521

$ws.Range("H42").Value = 34112.25
$ws.Range("J42").Value = 34112.25
$ws.Range("L42").Value = 34112.25
$ws.Range("N42").Value = -34868.25

$ws.Range("H43").Value = 15119.8
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 16399.75
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 16399.75
$ws.Range("M43").Value = -9851
$ws.Range("N43").Value = -16697.75

$ws.Range("H96").Value = 62501360
$ws.Range("I96").Value = 83334730
$ws.Range("J96").Value = 1250
$ws.Range("K96").Value = 83334730
$ws.Range("L96").Value = 1250
$ws.Range("M96").Value = -83333357
$ws.Range("N96").Value = -3996

$ws.Range("H105").Value = 43920
$ws.Range("J105").Value = 43920
$ws.Range("L105").Value = 43920
$ws.Range("N105").Value = -50908

$ws.Range("H122").Value = 1482.3667
$ws.Range("I122").Value = 1537.3462
$ws.Range("J122").Value = 1125
$ws.Range("K122").Value = 4612.0386
$ws.Range("L122").Value = 3375
$ws.Range("M122").Value = -2162.0386
$ws.Range("N122").Value = -8275

$ws.Range("H132").Value = 2161.8333
$ws.Range("I132").Value = 2228.8223
$ws.Range("K132").Value = 6686.466899999999
$ws.Range("M132").Value = -4156.466899999999
